# Update column G ("K") values on Sheet1, rows 2-31.
# These reflect a regenerated save of K (previously Strike#) per the
# commit message: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2 = 2
    3 = 4
    4 = 3
    5 = 8
    6 = 4
    7 = 4
    8 = 6
    9 = 0
    10 = 2
    11 = 3
    12 = 6
    13 = 7
    14 = 5
    15 = 3
    16 = 5
    17 = 10
    18 = 7
    19 = 4
    20 = 2
    21 = 6
    22 = 2
    23 = 4
    24 = 5
    25 = 4
    26 = 7
    27 = 2
    28 = 7
    29 = 2
    30 = 2
    31 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
